$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-11, columns A-G
# A: Colaborador_id, B: Colaborador_nome, C: Departamento,
# D: Motivo_da_ausência, E: Horas_de_ausência, F: Data_da_ausência, G: Salário

$data = @(
    @(56147, "Dra. Alícia Lopes",        "Financeiro",              "Consulta médica",    6, 45101, 6685.02),
    @(16692, "Davi Lucas da Cruz",       "Vendas",                  "Viagem de negócios", 2, 45089, 3205.28),
    @(75318, "Yasmin Silva",             "Operações",               "Viagem de negócios", 3, 45078, 12175.64),
    @(95005, "Rebeca Azevedo",           "P&D",                     "Outros",             1, 45106, 4088.04),
    @(71346, "João Pedro Pires",         "Atendimento ao Cliente",  "Viagem de negócios", 1, 45105, 7480.99),
    @(57799, "Leonardo Cardoso",         "TI",                      "Doença",             2, 45100, 5753.61),
    @(31382, "Elisa Monteiro",           "Marketing",               "Problemas pessoais", 5, 45086, 7433.63),
    @(84637, "Dr. Luiz Gustavo Porto",   "Vendas",                  "Doença",             7, 45097, 8137.08),
    @(56919, "Srta. Bárbara Peixoto",    "Engenharia",              "Outros",             1, 45099, 8471.28),
    @(47937, "Nathan Gomes",             "Vendas",                  "Problemas pessoais", 6, 45089, 8481.63)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}
